# Apply French translations to the Conditional Probability facilitator guide.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $new
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

Replace-Text "Video Title" "Titre de la vidéo"
Replace-Text "Topic" "Rubrique"
Replace-Text "Aim(s)" "Objectif(s)"
Replace-Text "Length" "Durée"
Replace-Text "Camp Location" "Lieu du camp"
Replace-Text "Facilitators" "Animateurs"
Replace-Text "N. of students" "N. des étudiants"
Replace-Text "Resources" "Les ressources"
Replace-Text "needed" "nécessaires"
Replace-Text "Preparations" "Préparations"
Replace-Text "Video time" "Temps de la vidéo"
Replace-Text "What facilitator does" "Ce que fait le facilitateur"
Replace-Text "What learners do" "Ce que font les apprenants"
Replace-Text "General VMC Video Introduction" "Vidéo générale introduisant le CVM"
Replace-Text "Video Introduction" "Video d'introduction"
Replace-Text "Facilitate the discussion: the “dealer” never told anything about the winning card, and nonetheless the “guesser” guessed right more times after the hint than before. Why? How can this be viewed in terms of the first experiment?" "Facilitate the discussion: the “dealer” never told anything about the winning card, and nonetheless the “guesser” guessed right more times after the hint than before. Pourquoi? How can this be viewed in terms of the first experiment?"
